# Refresh the cryptos price list (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) for most rows, and swaps the
# NEARProtocol/Aptos rows (36/37) which changed rank order.
#
# Note: several Price values look numeric (e.g. "582.49", "1.00",
# "0.0000117") but must stay plain text like the rest of the column, so
# a leading apostrophe forces Excel to store them as text instead of
# auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.490.63"
$ws.Range("E2").Value = "  +2.57%  "
$ws.Range("D3").Value = "3.483.25"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'582.49"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'147.75"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D7").Value = "3.482.49"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("E12").Value = "  +5.12%  "
$ws.Range("D13").Value = "4.075.98"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "'29.70"
$ws.Range("E14").Value = "  +5.26%  "
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").Value = "3.483.31"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "63.510.54"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("D20").Value = "'14.50"
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("D21").Value = "'9.35"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "'390.88"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'0.568"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").Value = "'74.93"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "3.621.68"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("E28").Value = "  -5.52%  "
$ws.Range("D29").Value = "'7.66"
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "'8.27"
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'1.37"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").Value = "'23.57"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.35"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'7.17"
$ws.Range("E37").Value = "  +2.93%  "
$ws.Range("D38").Value = "'1.61"
$ws.Range("E38").Value = "  +9.43%  "
$ws.Range("D39").Value = "'31.80"
$ws.Range("E39").Value = "  +12.04%  "
$ws.Range("D40").Value = "'169.68"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "3.520.45"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").Value = "'0.0766"
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("D43").Value = "'0.801"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("D45").Value = "'42.39"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("E46").Value = "  +3.40%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "2.613.06"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  +10.31%  "
$ws.Range("D50").Value = "'23.16"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").Value = "'6.82"
$ws.Range("E51").Value = "  +3.14%  "
